$wb = $excel.ActiveWorkbook

# --- Update data values ---
# addListItem sheet: "PuneAN" -> "PuneAO"
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "PuneAO"

# createUser sheet: 1091 -> 1092 (drives the CONCAT formulas in B2/F2)
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 1092

# --- Switch the active/selected tab from createUser to addListItem ---
$wsAddListItem.Activate()
$wsAddListItem.Select()
